$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 430.22223
$ws.Range("I6").Value = 430.22223
$ws.Range("K6").Value = 1290.66669
$ws.Range("M6").Value = -1178.66669
$ws.Range("H8").Value = 1303.5
$ws.Range("I8").Value = 61.142857
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 183.428571
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = -44.42857100000001
$ws.Range("N8").Value = -30278
$ws.Range("H96").Value = 1120.3334
$ws.Range("I96").Value = 1341.6
$ws.Range("J96").Value = 843.75
$ws.Range("K96").Value = 4024.8
$ws.Range("L96").Value = 2531.25
$ws.Range("M96").Value = -2651.8
$ws.Range("N96").Value = -5277.25
$ws.Range("H132").Value = 2875.111
$ws.Range("I132").Value = 2862.375
$ws.Range("J132").Value = 2977
$ws.Range("K132").Value = 8587.125
$ws.Range("L132").Value = 8931
$ws.Range("M132").Value = -6057.125
$ws.Range("N132").Value = -13991

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 12.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 12.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 317.5
$ws.Range("N26").ClearContents()
$ws.Range("H37").Value = 7258.5
$ws.Range("I37").Value = 4517
$ws.Range("K37").Value = 4517
$ws.Range("M37").Value = -4244
$ws.Range("H45").Value = 3914.1428
$ws.Range("I45").Value = 1199.5
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1199.5
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -822.5
$ws.Range("N45").Value = -5754
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 2479.2
$ws.Range("I74").Value = 1865.3334
$ws.Range("K74").Value = 1865.3334
$ws.Range("M74").Value = -991.3334
$ws.Range("H77").Value = 2479.2
$ws.Range("I77").Value = 1865.3334
$ws.Range("K77").Value = 9326.666999999999
$ws.Range("M77").Value = -4958.666999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1528.5
$ws.Range("I80").Value = 700
$ws.Range("K80").Value = 700
$ws.Range("M80").Value = 298
$ws.Range("H83").Value = 1528.5
$ws.Range("I83").Value = 700
$ws.Range("K83").Value = 3500
$ws.Range("M83").Value = 1492
$ws.Range("H105").Value = 1261.4286
$ws.Range("I105").Value = 1303.1666
$ws.Range("K105").Value = 1303.1666
$ws.Range("M105").Value = 443.8334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1503
$ws.Range("I31").Value = 1503
$ws.Range("K31").Value = 1503
$ws.Range("M31").Value = -1208
$ws.Range("H34").Value = 1503
$ws.Range("I34").Value = 1503
$ws.Range("K34").Value = 1503
$ws.Range("M34").Value = -1301
$ws.Range("H60").Value = 21000
$ws.Range("I60").Value = 22500
$ws.Range("K60").Value = 22500
$ws.Range("M60").Value = -21989
$ws.Range("H107").Value = 1503.2
$ws.Range("I107").Value = 1564.6666
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 1564.6666
$ws.Range("L107").Value = 950
$ws.Range("M107").Value = 355.3334
$ws.Range("N107").Value = -4790
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22080
$ws.Range("H132").Value = 3497.5
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -20045

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5418.8335
$ws.Range("I11").Value = 4.6666665
$ws.Range("J11").Value = 10833
$ws.Range("K11").Value = 13.9999995
$ws.Range("L11").Value = 32499
$ws.Range("M11").Value = 126.0000005
$ws.Range("N11").Value = -32779
$ws.Range("H21").Value = 850
$ws.Range("I21").Value = 750
$ws.Range("K21").Value = 2250
$ws.Range("M21").Value = -2077
$ws.Range("H25").Value = 491.66666
$ws.Range("I25").Value = 400
$ws.Range("J25").Value = 675
$ws.Range("K25").Value = 1200
$ws.Range("L25").Value = 2025
$ws.Range("M25").Value = -1031
$ws.Range("N25").Value = -2363
$ws.Range("H30").Value = 491.66666
$ws.Range("I30").Value = 400
$ws.Range("J30").Value = 675
$ws.Range("K30").Value = 1200
$ws.Range("L30").Value = 2025
$ws.Range("M30").Value = -1098
$ws.Range("N30").Value = -2229
$ws.Range("H38").Value = 54.636364
$ws.Range("I38").Value = 34
$ws.Range("K38").Value = 102
$ws.Range("M38").Value = 245
$ws.Range("H40").Value = 72
$ws.Range("I40").Value = 60.57143
$ws.Range("J40").Value = 83.42856999999999
$ws.Range("K40").Value = 242.28572
$ws.Range("L40").Value = 333.71428
$ws.Range("M40").Value = -173.28572
$ws.Range("N40").Value = -471.71428
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2510000
$ws.Range("I3").Value = 2510000
$ws.Range("K3").Value = 2510000
$ws.Range("M3").Value = -2509884
$ws.Range("H5").Value = 37500
$ws.Range("I5").Value = 25000
$ws.Range("K5").Value = 25000
$ws.Range("M5").Value = -24888
$ws.Range("H132").Value = 3464.6875
$ws.Range("I132").Value = 3429
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 10287
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7757
$ws.Range("N132").Value = -17060
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3616.0881
$ws.Range("I46").Value = 3562.9565
$ws.Range("K46").Value = 3562.9565
$ws.Range("M46").Value = -3374.9565
$ws.Range("H132").Value = 2333
$ws.Range("J132").Value = 2499.5
$ws.Range("L132").Value = 7498.5
$ws.Range("N132").Value = -12558.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000.5
$ws.Range("J64").Value = 30000.5
$ws.Range("L64").Value = 30000.5
$ws.Range("N64").Value = -30496.5
$ws.Range("H67").Value = 30000.5
$ws.Range("J67").Value = 30000.5
$ws.Range("L67").Value = 30000.5
$ws.Range("N67").Value = -31716.5
$ws.Range("H81").Value = 468.8
$ws.Range("I81").Value = 468.8
$ws.Range("K81").Value = 937.6
$ws.Range("M81").Value = 123.4
$ws.Range("H84").Value = 468.8
$ws.Range("I84").Value = 468.8
$ws.Range("K84").Value = 4688
$ws.Range("M84").Value = 616
$ws.Range("H100").Value = 3031608.8
$ws.Range("I100").Value = 5362306
$ws.Range("J100").Value = 1702.9
$ws.Range("K100").Value = 10724612
$ws.Range("L100").Value = 3405.8
$ws.Range("M100").Value = -10724071
$ws.Range("N100").Value = -4487.8
